$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D to make room for "Status"
$ws.Columns("D").Insert()

# Update header row
$ws.Range("D1").Value = "Status"
$ws.Range("G1").Value = "Oct_2025"

# Update ISIN / Stock Name for rows 15 and 16 (swap of holdings)
$ws.Range("A15").Value = "INE326A01037"
$ws.Range("B15").Value = "Lupin Limited"
$ws.Range("A16").Value = "INE044A01036"
$ws.Range("B16").Value = "Sun Pharmaceutical Industries Limited"

# Fill in Status column and refreshed data values for each data row
$ws.Range("D2").Value = "Adding Consistently"
$ws.Range("E2").Value = 8.831271
$ws.Range("F2").Value = 8.50142
$ws.Range("G2").Value = 8.25627
$ws.Range("H2").Value = 0.3298509999999997
$ws.Range("I2").Value = 0.5750009999999985

$ws.Range("D3").Value = "Adding Consistently"
$ws.Range("E3").Value = 8.447489
$ws.Range("F3").Value = 8.051306
$ws.Range("G3").Value = 7.965679
$ws.Range("H3").Value = 0.3961829999999988
$ws.Range("I3").Value = 0.4818099999999994

$ws.Range("D4").Value = "Reducing"
$ws.Range("E4").Value = 8.441843
$ws.Range("F4").Value = 9.582952
$ws.Range("G4").Value = 6.526537
$ws.Range("H4").Value = -1.141109
$ws.Range("I4").Value = 1.915306

$ws.Range("D5").Value = "Adding Consistently"
$ws.Range("E5").Value = 7.705935
$ws.Range("F5").Value = 7.094618
$ws.Range("G5").Value = 6.445576
$ws.Range("H5").Value = 0.6113170000000006
$ws.Range("I5").Value = 1.260359

$ws.Range("D6").Value = "Reducing Consistently"
$ws.Range("E6").Value = 7.241881
$ws.Range("F6").Value = 7.604171
$ws.Range("G6").Value = 7.887981
$ws.Range("H6").Value = -0.3622899999999998
$ws.Range("I6").Value = -0.6460999999999997

$ws.Range("D7").Value = "Adding Consistently"
$ws.Range("E7").Value = 6.950458
$ws.Range("F7").Value = 6.230534
$ws.Range("G7").Value = 5.801788
$ws.Range("H7").Value = 0.7199240000000007
$ws.Range("I7").Value = 1.14867

$ws.Range("D8").Value = "Reducing Consistently"
$ws.Range("E8").Value = 6.425974
$ws.Range("F8").Value = 6.601
$ws.Range("G8").Value = 6.505723
$ws.Range("H8").Value = -0.1750259999999999
$ws.Range("I8").Value = -0.07974899999999963

$ws.Range("D9").Value = "Reducing Consistently"
$ws.Range("E9").Value = 5.869715
$ws.Range("F9").Value = 5.997429
$ws.Range("G9").Value = 6.16882
$ws.Range("H9").Value = -0.1277140000000001
$ws.Range("I9").Value = -0.299105

$ws.Range("D10").Value = "Reducing Consistently"
$ws.Range("E10").Value = 5.635718
$ws.Range("F10").Value = 5.871914
$ws.Range("G10").Value = 6.138305
$ws.Range("H10").Value = -0.2361960000000005
$ws.Range("I10").Value = -0.5025870000000001

$ws.Range("D11").Value = "Adding Consistently"
$ws.Range("E11").Value = 3.628896
$ws.Range("F11").Value = 3.603264
$ws.Range("G11").Value = 3.584167
$ws.Range("H11").Value = 0.02563200000000032
$ws.Range("I11").Value = 0.04472900000000024

$ws.Range("D12").Value = "Fresh Entry"
$ws.Range("E12").Value = 3.520651
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 3.520651
$ws.Range("I12").Value = 3.520651

$ws.Range("D13").Value = "Reducing Consistently"
$ws.Range("E13").Value = 2.880442
$ws.Range("F13").Value = 3.099855
$ws.Range("G13").Value = 3.138804
$ws.Range("H13").Value = -0.2194129999999999
$ws.Range("I13").Value = -0.258362

$ws.Range("D14").Value = "Adding Consistently"
$ws.Range("E14").Value = 1.52235
$ws.Range("F14").Value = 1.334298
$ws.Range("G14").Value = 1.410462
$ws.Range("H14").Value = 0.1880520000000001
$ws.Range("I14").Value = 0.111888

$ws.Range("D15").Value = "Complete Exit"
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 2.913663
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = -2.913663

$ws.Range("D16").Value = "Complete Exit"
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 3.401542
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = -3.401542

$ws.Range("D17").Value = "Complete Exit"
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 2.801618
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = -2.801618

$ws.Range("D18").Value = "Complete Exit"
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 4.091851
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = -4.091851
$ws.Range("I18").Value = 0

$ws.Range("D19").Value = "Complete Exit"
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 10.010286
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = -10.010286
$ws.Range("I19").Value = 0
